# Slide 7 ("... pore Commit changes e click korle ...") text-box edit:
#  1. Merge the two single-space runs (between the 2nd "ক্লিক" and "করে")
#     into a single run containing two spaces.
#  2. Extend the closing "।" with a trailing space and append a brand new
#     sentence describing the "Commit changes" step (mostly bold).
#  3. Start a fresh, empty trailing paragraph.
#  4. Grow the text box to its new autofit height.

$p  = $ppt.ActivePresentation
$s  = $p.Slides.Item(7)
$sh = $s.Shapes.Item(1)
$tr = $sh.TextFrame.TextRange

# --- 1. Merge the double-space runs (1-based char index 236, length 2) ----
$dbl = $tr.Characters(236, 2)
$dbl.Text = "  "

# --- 2. Turn the trailing "।" into "। " -----------------------------------
$closer = $tr.Characters($tr.Text.Length, 1)
$closer.Text = "। "

# --- 3. Append the new sentence, then re-run the formatting over it so it
#        is split into the same runs a real PowerPoint edit session would
#        leave behind (plain lead-in, bold "Commit changes ..." tail). ----
$newSentence = "এর পর Commit changes এ ক্লিক করলে File Save হয়ে যাবে।"
$tr.InsertAfter($newSentence)

$segments = @(
  @{ Start = 282; Length = 2;  Bold = $false },  # এর
  @{ Start = 284; Length = 1;  Bold = $false },  # " "
  @{ Start = 285; Length = 2;  Bold = $false },  # পর
  @{ Start = 287; Length = 1;  Bold = $false },  # " "
  @{ Start = 288; Length = 7;  Bold = $true  },  # "Commit "
  @{ Start = 295; Length = 10; Bold = $true  },  # "changes এ "
  @{ Start = 305; Length = 5;  Bold = $true  },  # ক্লিক
  @{ Start = 310; Length = 1;  Bold = $true  },  # " "
  @{ Start = 311; Length = 4;  Bold = $true  },  # করলে
  @{ Start = 315; Length = 11; Bold = $true  },  # " File Save "
  @{ Start = 326; Length = 4;  Bold = $true  },  # হয়ে
  @{ Start = 330; Length = 1;  Bold = $true  },  # " "
  @{ Start = 331; Length = 4;  Bold = $true  },  # যাবে
  @{ Start = 335; Length = 1;  Bold = $true  }   # ।
)

foreach ($seg in $segments) {
    $rng = $tr.Characters($seg.Start, $seg.Length)
    if ($seg.Bold) {
        # Explicit bold toggle -> run gets b="1".
        $rng.Font.Bold = $true
    } else {
        # Re-assert the same text to force a run boundary here without
        # stamping an explicit b="0" on a run that was never bold.
        $rng.Text = $rng.Text
    }
}

# --- 4. New empty trailing paragraph ---------------------------------------
$tr.InsertAfter([char]13)

# --- 5. Resize the text box for the new (taller) autofit content ----------
$sh.Height = 5016758 / 914400 * 72
